$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 57: Mecca, Saudi Arabia (Great Mosque) ---
$ws.Range("F57").Value = "G0PC9JDC-2Y"
$ws.Range("E57").Value = "Saudi Arabia"
$ws.Range("B57").Value = "21.42254924107721, 39.82620067445773"
$ws.Range("D57").Value = "Mecca"
$ws.Range("C57").Value = "Great Mosque - Live broadcast || The Holy Quran Channel"
$ws.Range("A57").Value = "LIVE, RELIGION"

# --- Row 58: Medina, Saudi Arabia (Sunnah Channel) ---
$ws.Range("F58").Value = "Kt7hKHlArl8"
$ws.Range("D58").Value = "Medina"
$ws.Range("C58").Value = "Live broadcast || Sunnah Channel"
$ws.Range("B58").Value = "24.46893447433643, 39.61111417777757"
$ws.Range("E58").Value = "Saudi Arabia"
$ws.Range("A58").Value = "LIVE, RELIGION, CITY"

# Copy the bordered style used by the rest of the table (columns A, C, D, E)
# from the last existing data row (56) down onto the two new rows, leaving
# B and F (which carry no explicit style) untouched.
$ws.Range("A56").Copy()
$ws.Range("A57").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A58").PasteSpecial(-4122)

$ws.Range("C56").Copy()
$ws.Range("C57").PasteSpecial(-4122)
$ws.Range("C58").PasteSpecial(-4122)

$ws.Range("D56").Copy()
$ws.Range("D57").PasteSpecial(-4122)
$ws.Range("D58").PasteSpecial(-4122)

$ws.Range("E56").Copy()
$ws.Range("E57").PasteSpecial(-4122)
$ws.Range("E58").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Mirror the view state left behind after editing (selection resting on
# the new last row).
$ws.Range("A59").Select() | Out-Null
